$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Seed the new shared strings in the exact order they appear in the
#    target sharedStrings.xml (Sam=15, Jess=16, Part-Time=17, No=18, Yes=19,
#    Git=20). Writing string literals to cells in this order makes the
#    unique-string table grow with the same ordering as the target file.
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "Sam"
$ws.Range("E1").Value = "Jess"
$ws.Range("A28").Value = "Part-Time"
$ws.Range("B28").Value = "No"
$ws.Range("D28").Value = "Yes"
$ws.Range("A25").Value = "Git"

# ---------------------------------------------------------------------------
# 2. Table 1 (weekly hours, rows 1-15): add Jess/Sam header + new "Jess"
#    column (E) values, then repoint the Average/SUM helper formulas at
#    A:E instead of B:F.
# ---------------------------------------------------------------------------
$ws.Range("F17").Value = "Sam"
$ws.Range("E17").Value = "Jess"

$tbl1 = @{
    2  = 2
    3  = 4
    4  = 3
    5  = 3
    6  = 5
    7  = 5
    8  = 5
    9  = 5
    10 = 10
    11 = 5
    12 = 5
    13 = 3
}
foreach ($r in $tbl1.Keys) {
    $ws.Range("E${r}").Value = $tbl1[$r]
}

# Row 2 uses standalone formulas; rows 3-13 are a shared-formula block, so
# re-writing the anchor (row 3) and then every other row keeps the engine's
# shared-formula grouping while updating each row's own relative reference.
$ws.Range("I2").Formula = "=AVERAGE(A2:E2)"
$ws.Range("J2").Formula = "=SUM(A2:E2)"
for ($r = 3; $r -le 13; $r++) {
    $ws.Range("I${r}").Formula = "=AVERAGE(A${r}:E${r})"
    $ws.Range("J${r}").Formula = "=SUM(A${r}:E${r})"
}

# Row 15 totals for the two new columns.
$ws.Range("E15").Formula = "=SUM(E2:E13)"
$ws.Range("F15").Formula = "=SUM(F2:F13)"

# J-column "traffic light" cell styles that move around because the SUM
# total for several weeks changed once Jess's hours were added in.
$jStyles = @{
    2  = "Bad"
    3  = "Bad"
    4  = "Bad"
    5  = "Neutral"
    6  = "Bad"
    7  = "Neutral"
    8  = "Neutral"
    9  = "Neutral"
    10 = "Good"
    11 = "Neutral"
    12 = "Good"
    13 = "Good"
}
foreach ($r in $jStyles.Keys) {
    $ws.Range("J${r}").Style = $jStyles[$r]
}

# ---------------------------------------------------------------------------
# 3. Table 2 (skills, rows 17-25): add the "Jess" column (E), a new "Git"
#    row (25), and widen the Average formula to include column E.
# ---------------------------------------------------------------------------
$tbl2 = @{
    18 = 1
    19 = 8
    20 = 5
    21 = 7
    22 = 8
    23 = 2
    24 = 1
}
foreach ($r in $tbl2.Keys) {
    $ws.Range("E${r}").Value = $tbl2[$r]
}

# Jess turns out to have spent a random evening doing some Java that week.
$ws.Range("B21").Value = 3

for ($r = 18; $r -le 24; $r++) {
    $ws.Range("I${r}").Formula = "=AVERAGE(B${r}:F${r})"
}

# New "Git" skill row.
$ws.Range("B25").Value = 7
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 8
$ws.Range("I25").Formula = "=AVERAGE(B25:E25)"

# ---------------------------------------------------------------------------
# 4. New "Part-Time" lookup row (28).
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = "No"
$ws.Range("E28").Value = "Yes"
$ws.Range("F28").Value = "No"

# ---------------------------------------------------------------------------
# 5. Column width for A (auto-fit-ish, new longer labels like "Part-Time").
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 10

# ---------------------------------------------------------------------------
# 6. Charts: extend/append the "Jess" series on all three charts, and widen
#    the skills charts' categories to include the new "Git" row.
# ---------------------------------------------------------------------------
$chartObjs = $ws.ChartObjects()

# Chart "Chart 3" (first chart object) and "Chart 2" (second) both plot the
# weekly table (rows 2-13); add a 4th series for Jess (column E).
for ($i = 1; $i -le 2; $i++) {
    $chart = $chartObjs.Item($i).Chart
    $ser = $chart.SeriesCollection().NewSeries()
    $ser.Name = "=Sheet1!`$E`$1"
    $ser.Values = "=Sheet1!`$E`$2:`$E`$13"
}

# Chart "Chart 1" (third chart object) plots the skills table (rows 18-24);
# widen the existing three series to include the new Git row, then append
# the Jess series across rows 18-25.
$skillsChart = $chartObjs.Item(3).Chart
for ($i = 1; $i -le 3; $i++) {
    $ser = $skillsChart.SeriesCollection().Item($i)
    $col = [char](66 + $i - 1)
    $ser.XValues = "=Sheet1!`$A`$18:`$A`$25"
    $ser.Values = "=Sheet1!`$${col}`$18:`$${col}`$25"
}
$jessSkills = $skillsChart.SeriesCollection().NewSeries()
$jessSkills.Name = "=Sheet1!`$E`$17"
$jessSkills.XValues = "=Sheet1!`$A`$18:`$A`$25"
$jessSkills.Values = "=Sheet1!`$E`$18:`$E`$25"

# ---------------------------------------------------------------------------
# 7. Move the skills chart down/right slightly (matches its new anchor) and
#    refresh the current selection.
# ---------------------------------------------------------------------------
$skillsChartObj = $chartObjs.Item(3)
$skillsChartObj.Left = $skillsChartObj.Left + 7
$skillsChartObj.Top = $skillsChartObj.Top - 2

$ws.Range("H26").Select()
